$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to reflect new report column meanings
$ws.Range("C3").Value = "Hours"
$ws.Range("B5").Value = "Total"

# Update the active selection/cell
$ws.Range("E8").Select()
